$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-09-01 22:49:47"

$wsZhCn.Range("H3").Value = "2016-09-01 22:49:43"
$wsZhCn.Range("K3").Value = "2016-09-01 22:50:09"

$wsDeDe.Range("H3").Value = "2016-09-01 22:49:47"
$wsDeDe.Range("K3").Value = "2016-09-01 22:50:20"
